# Update the "Spleen" row (row 16):
#  - rename the Variant/Detail text
#  - add the new YouTube clip link with a hyperlink (matching the style of
#    the other hyperlink cells D3/D8/D13)
# and leave the selection where the user ended up after the edit (D19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = "Splenic calcification with posterior shadowing"
$ws.Range("D16").Value = "https://youtu.be/qushjTAy6XQ "
$ws.Hyperlinks.Add($ws.Range("D16"), "https://youtu.be/qushjTAy6XQ")
$ws.Range("D16").Style = "Collegamento ipertestuale"

$ws.Range("D19").Select()
